# Adapt column header formatting to respective input file names (#7)
#  - rename "<header>_old" -> "<header>_FV2310"  (columns A-J)
#  - rename "<header>_new" -> "<header>_FV2404"  (columns L-U)
#  - turn the used range A1:U69 into a real Excel Table ("Table1")
#  - freeze the header row (top row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (suffix swap) -------------------------------
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old  = [string]$cell.Value()
    $cell.Value = $old.Replace("_old", "_FV2310")
}

for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old  = [string]$cell.Value()
    $cell.Value = $old.Replace("_new", "_FV2404")
}

# --- 2. Convert the data range into a real Excel Table -------------------
# The header row (row 1) already carries bold/filled formatting. Creating
# the ListObject directly on top of that range makes Excel capture the
# existing look into a header-row dxf override. Building the table on a
# bare scratch range first (no special formatting) and then resizing it
# onto the real data avoids that extra, unwanted style override.
$lastRow = 101
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($lastRow, $c).Value = "c$c"
    $ws.Cells.Item($lastRow + 1, $c).Value = "x"
}
$scratch = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow + 1, $lastCol))
$tbl = $ws.ListObjects.Add(1, $scratch, 0, 1)
$tbl.Name = "Table1"

$target = $ws.Range("A1:U69")
$tbl.Resize($target)
$tbl.TableStyle = ""

# Re-push the (already-correct) header text so the table's column
# collection picks the real names up instead of the scratch placeholders.
for ($c = 1; $c -le $lastCol; $c++) {
    $headerValue = [string]$ws.Cells.Item(1, $c).Value()
    $ws.Cells.Item(1, $c).Value = $headerValue
}

# Drop the scratch rows again so they don't linger outside the table.
$cleanup = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow + 1, $lastCol))
$cleanup.Clear()

# --- 3. Freeze the header row --------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
